## Handleiding applicatie bewoners - content edit
## Applies the following user-visible changes:
##   1. Fix "wif|i" typo-split text: merge the wifi sentence into a single
##      run (removing the stray leftover _GoBack bookmark that had been
##      sitting in the middle of the word "wifi").
##   2. "... Dit kunt u doen door rechts onderin het scherm ..." -> "... links ..."
##      (how to reach the quick-contacts list).
##   3. "drukt u 3 seconden ... Zodra deze 3 seconden ..." -> both "3"s become "2"
##      for the noodknop (emergency button) instructions, with Word's
##      "_GoBack" last-edit bookmark left at the second edit location.
##   4. The cached page-number field result in the default footer updates
##      from "3" to "4".

$d = $word.ActiveDocument

function Replace-InContext {
    param(
        $Doc,
        [string]$Context,
        [int]$Offset,
        [int]$Length,
        [string]$NewText
    )
    $probe = $Doc.Content
    $found = $probe.Find.Execute($Context, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $null
    }
    $start = $probe.Start + $Offset
    $end = $start + $Length
    $target = $Doc.Range($start, $end)
    $target.Text = $NewText
    return $target
}

## 1. "Voor begeleiding met het verbinden van het wifi-netwerk, raadpleeg de
##    handleiding verbinden met wif" + [_GoBack] + "i. Deze kunt u vinden in
##    de handleiding bundel." -> one consolidated run (bookmark disappears).
$d.Content.Find.Execute(
    "raadpleeg de handleiding verbinden met wifi. Deze kunt u vinden in de handleiding bundel.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "raadpleeg de handleiding verbinden met wifi. Deze kunt u vinden in de handleiding bundel.",
    2) | Out-Null

## 2. "... doen door rechts onderin het scherm op de groene belknop ..." -> "links"
Replace-InContext $d "doen door rechts onderin" 10 6 "links" | Out-Null

## 3. Noodknop paragraph: both occurrences of "3 seconden" -> "2 seconden".
Replace-InContext $d "drukt u 3 seconden" 8 1 "2" | Out-Null
$secondEdit = Replace-InContext $d "Zodra deze 3 seconden" 11 1 "2"

## Word drops its "_GoBack" bookmark at the last place text was edited;
## re-create it (zero length, right after the second "2") since the earlier
## wifi-sentence edit already consumed/removed the original one.
if ($secondEdit -ne $null) {
    $goBackPos = $secondEdit.End
    $goBackRange = $d.Range($goBackPos, $goBackPos)
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
}

## 4. Default footer's cached PAGE field result: "3" -> "4".
$footerStory = $d.StoryRanges(9)
$footerStory.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "4", 2) | Out-Null
